$d = $word.ActiveDocument

# 1. "aki" -> "akí" in the phrase about "e komportashon aki sosodé"
$d.Content.Find.Execute(
    "e komportashon aki sosodé",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "e komportashon akí sosodé",
    2)

# 2. "enfoká" -> "enfokando" in "e lo ta enfoká riba e yu"
$d.Content.Find.Execute(
    "lo ta enfoká riba e yu ku el a skohe",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "lo ta enfokando riba e yu ku el a skohe",
    2)

# 3. insert "e " -> "ku e ta siña den e programa"
$d.Content.Find.Execute(
    "ku ta siña den e programa",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "ku e ta siña den e programa",
    2)

# 4. "aplikahson" -> "aplikashon"
$d.Content.Find.Execute(
    "kon pa baha e aplikahson si nan",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "kon pa baha e aplikashon si nan",
    2)
